# Auto-generated: apply numeric cell updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 454.7
$ws.Range("J2").Value = 618.75
$ws.Range("L2").Value = 618.75
$ws.Range("N2").Value = -844.75
$ws.Range("H11").Value = 299.5
$ws.Range("I11").Value = 299.5
$ws.Range("K11").Value = 299.5
$ws.Range("M11").Value = -159.5
$ws.Range("H12").Value = 1278.25
$ws.Range("I12").Value = 977.8
$ws.Range("J12").Value = 1779
$ws.Range("K12").Value = 977.8
$ws.Range("L12").Value = 1779
$ws.Range("M12").Value = -807.8
$ws.Range("N12").Value = -2119
$ws.Range("H64").Value = 4600
$ws.Range("I64").Value = 6000
$ws.Range("K64").Value = 6000
$ws.Range("M64").Value = -5752
$ws.Range("H67").Value = 4600
$ws.Range("I67").Value = 6000
$ws.Range("K67").Value = 6000
$ws.Range("M67").Value = -5142
$ws.Range("H76").Value = 6401
$ws.Range("I76").Value = 5365.375
$ws.Range("J76").Value = 7229.5
$ws.Range("K76").Value = 5365.375
$ws.Range("L76").Value = 7229.5
$ws.Range("M76").Value = -5050.375
$ws.Range("N76").Value = -7859.5
$ws.Range("H79").Value = 6401
$ws.Range("I79").Value = 5365.375
$ws.Range("J79").Value = 7229.5
$ws.Range("K79").Value = 5365.375
$ws.Range("L79").Value = 7229.5
$ws.Range("M79").Value = -4273.375
$ws.Range("N79").Value = -9413.5
$ws.Range("H86").Value = 4072.2354
$ws.Range("I86").Value = 1873
$ws.Range("J86").Value = 4748.923
$ws.Range("K86").Value = 1873
$ws.Range("L86").Value = 4748.923
$ws.Range("M86").Value = -750
$ws.Range("N86").Value = -6994.923
$ws.Range("H89").Value = 4072.2354
$ws.Range("I89").Value = 1873
$ws.Range("J89").Value = 4748.923
$ws.Range("K89").Value = 9365
$ws.Range("L89").Value = 23744.615
$ws.Range("M89").Value = -3749
$ws.Range("N89").Value = -34976.615
$ws.Range("H116").Value = 8652.143
$ws.Range("I116").Value = 8822.75
$ws.Range("J116").Value = 8424.666999999999
$ws.Range("K116").Value = 8822.75
$ws.Range("L116").Value = 8424.666999999999
$ws.Range("M116").Value = -5380.75
$ws.Range("N116").Value = -15308.667
$ws.Range("H118").Value = 2607.5
$ws.Range("I118").Value = 2607.5
$ws.Range("K118").Value = 7822.5
$ws.Range("M118").Value = -6165.5
$ws.Range("H125").Value = 2766.6667
$ws.Range("I125").Value = 2900
$ws.Range("K125").Value = 26100
$ws.Range("M125").Value = -23640
$ws.Range("H129").Value = 3897.3333
$ws.Range("I129").Value = 1032.2
$ws.Range("J129").Value = 4999.3076
$ws.Range("K129").Value = 3096.6
$ws.Range("L129").Value = 14997.9228
$ws.Range("M129").Value = 1903.4
$ws.Range("N129").Value = -24997.9228
$ws.Range("H131").Value = 4321.5625
$ws.Range("I131").Value = 914.6
$ws.Range("K131").Value = 2743.8
$ws.Range("M131").Value = 2296.2
$ws.Range("H132").Value = 1903.9546
$ws.Range("I132").Value = 1625.7368
$ws.Range("K132").Value = 4877.2104
$ws.Range("M132").Value = -2347.2104
$ws.Range("H135").Value = 1132.9375
$ws.Range("I135").Value = 552.4167
$ws.Range("K135").Value = 4971.7503
$ws.Range("M135").Value = -2436.7503
$ws.Range("H137").Value = 2537.9092
$ws.Range("I137").Value = 1391.091
$ws.Range("J137").Value = 3684.7273
$ws.Range("K137").Value = 4173.272999999999
$ws.Range("L137").Value = 11054.1819
$ws.Range("M137").Value = -1623.272999999999
$ws.Range("N137").Value = -16154.1819
$ws.Range("H138").Value = 3760.985
$ws.Range("I138").Value = 1222.7273
$ws.Range("J138").Value = 4259.5713
$ws.Range("K138").Value = 3668.1819
$ws.Range("L138").Value = 12778.7139
$ws.Range("M138").Value = 1471.8181
$ws.Range("N138").Value = -23058.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1401.5883
$ws.Range("I61").Value = 1401.5883
$ws.Range("K61").Value = 1401.5883
$ws.Range("M61").Value = -1189.5883
$ws.Range("H63").Value = 5831.5
$ws.Range("I63").Value = 5107.778
$ws.Range("J63").Value = 6555.222
$ws.Range("K63").Value = 5107.778
$ws.Range("L63").Value = 6555.222
$ws.Range("M63").Value = -4421.778
$ws.Range("N63").Value = -7927.222
$ws.Range("H66").Value = 5831.5
$ws.Range("I66").Value = 5107.778
$ws.Range("J66").Value = 6555.222
$ws.Range("K66").Value = 25538.89
$ws.Range("L66").Value = 32776.11
$ws.Range("M66").Value = -22106.89
$ws.Range("N66").Value = -39640.11
$ws.Range("H74").Value = 4153.278
$ws.Range("I74").Value = 1253
$ws.Range("K74").Value = 1253
$ws.Range("M74").Value = -379
$ws.Range("H77").Value = 4153.278
$ws.Range("I77").Value = 1253
$ws.Range("K77").Value = 6265
$ws.Range("M77").Value = -1897
$ws.Range("H97").Value = 529.4167
$ws.Range("I97").Value = 491.18182
$ws.Range("K97").Value = 491.18182
$ws.Range("M97").Value = 4.818179999999984
$ws.Range("H102").Value = 1245.1666
$ws.Range("I102").Value = 894.5
$ws.Range("K102").Value = 894.5
$ws.Range("M102").Value = 727.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 717463.2
$ws.Range("I122").Value = 1252085.6
$ws.Range("K122").Value = 3756256.8
$ws.Range("M122").Value = -3753806.8
$ws.Range("H132").Value = 1072.36
$ws.Range("I132").Value = 1078.6957
$ws.Range("K132").Value = 3236.0871
$ws.Range("M132").Value = -706.0870999999997
$ws.Range("H136").Value = 1401.5883
$ws.Range("I136").Value = 1401.5883
$ws.Range("K136").Value = 4204.7649
$ws.Range("M136").Value = -1654.7649

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 529.2222
$ws.Range("J22").Value = 425
$ws.Range("L22").Value = 425
$ws.Range("N22").Value = -771
$ws.Range("H99").Value = 549.1429000000001
$ws.Range("I99").Value = 432.33334
$ws.Range("K99").Value = 432.33334
$ws.Range("M99").Value = 1065.66666
$ws.Range("H134").Value = 3247.6365
$ws.Range("I134").Value = 3164.238
$ws.Range("K134").Value = 9492.714
$ws.Range("M134").Value = -6957.714
$ws.Range("H141").Value = 200000
$ws.Range("J141").Value = 200000
$ws.Range("L141").Value = 200000
$ws.Range("N141").Value = -210360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 540.2
$ws.Range("I16").Value = 672.8
$ws.Range("J16").Value = 496
$ws.Range("K16").Value = 672.8
$ws.Range("L16").Value = 496
$ws.Range("M16").Value = -385.8
$ws.Range("N16").Value = -1070
$ws.Range("H31").Value = 4253.6
$ws.Range("I31").Value = 2053.5
$ws.Range("J31").Value = 6768
$ws.Range("K31").Value = 2053.5
$ws.Range("L31").Value = 6768
$ws.Range("M31").Value = -1758.5
$ws.Range("N31").Value = -7358
$ws.Range("H34").Value = 4253.6
$ws.Range("I34").Value = 2053.5
$ws.Range("J34").Value = 6768
$ws.Range("K34").Value = 2053.5
$ws.Range("L34").Value = 6768
$ws.Range("M34").Value = -1851.5
$ws.Range("N34").Value = -7172
$ws.Range("H35").Value = 5261.25
$ws.Range("I35").Value = 5261.25
$ws.Range("K35").Value = 5261.25
$ws.Range("M35").Value = -4967.25
$ws.Range("H58").Value = 3151.05
$ws.Range("I58").Value = 1626.8
$ws.Range("J58").Value = 4675.3
$ws.Range("K58").Value = 1626.8
$ws.Range("L58").Value = 4675.3
$ws.Range("M58").Value = -1423.8
$ws.Range("N58").Value = -5081.3
$ws.Range("H62").Value = 83519.60000000001
$ws.Range("I62").Value = 4199.6665
$ws.Range("J62").Value = 202499.5
$ws.Range("K62").Value = 4199.6665
$ws.Range("L62").Value = 202499.5
$ws.Range("M62").Value = -3575.6665
$ws.Range("N62").Value = -203747.5
$ws.Range("H65").Value = 83519.60000000001
$ws.Range("I65").Value = 4199.6665
$ws.Range("J65").Value = 202499.5
$ws.Range("K65").Value = 20998.3325
$ws.Range("L65").Value = 1012497.5
$ws.Range("M65").Value = -17878.3325
$ws.Range("N65").Value = -1018737.5
$ws.Range("H86").Value = 14248.667
$ws.Range("J86").Value = 16499.25
$ws.Range("L86").Value = 16499.25
$ws.Range("N86").Value = -18745.25
$ws.Range("H89").Value = 14248.667
$ws.Range("J89").Value = 16499.25
$ws.Range("L89").Value = 82496.25
$ws.Range("N89").Value = -93728.25
$ws.Range("H107").Value = 153
$ws.Range("I107").Value = 153
$ws.Range("K107").Value = 153
$ws.Range("M107").Value = 1767
$ws.Range("H113").Value = 540.2
$ws.Range("I113").Value = 672.8
$ws.Range("J113").Value = 496
$ws.Range("K113").Value = 672.8
$ws.Range("L113").Value = 496
$ws.Range("M113").Value = 1497.2
$ws.Range("N113").Value = -4836
$ws.Range("H122").Value = 7510.3184
$ws.Range("I122").Value = 7314.5454
$ws.Range("J122").Value = 7706.091
$ws.Range("K122").Value = 21943.6362
$ws.Range("L122").Value = 23118.273
$ws.Range("M122").Value = -19493.6362
$ws.Range("N122").Value = -28018.273
$ws.Range("H132").Value = 2598.3215
$ws.Range("J132").Value = 2941.0908
$ws.Range("L132").Value = 8823.2724
$ws.Range("N132").Value = -13883.2724
$ws.Range("H134").Value = 4226.0713
$ws.Range("I134").Value = 3451.375
$ws.Range("J134").Value = 5259
$ws.Range("K134").Value = 10354.125
$ws.Range("L134").Value = 15777
$ws.Range("M134").Value = -7819.125
$ws.Range("N134").Value = -20847
$ws.Range("H136").Value = 3151.05
$ws.Range("I136").Value = 1626.8
$ws.Range("J136").Value = 4675.3
$ws.Range("K136").Value = 4880.4
$ws.Range("L136").Value = 14025.9
$ws.Range("M136").Value = -2330.4
$ws.Range("N136").Value = -19125.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125242.125
$ws.Range("I2").Value = 143076.72
$ws.Range("K2").Value = 858460.3200000001
$ws.Range("M2").Value = -858347.3200000001
$ws.Range("H26").Value = 896.25
$ws.Range("I26").Value = 195
$ws.Range("K26").Value = 585
$ws.Range("M26").Value = -297
$ws.Range("H39").Value = 2121.818
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H55").Value = 127537.5
$ws.Range("J55").Value = 4725
$ws.Range("L55").Value = 14175
$ws.Range("N55").Value = -14529
$ws.Range("H92").Value = 2510
$ws.Range("I92").Value = 2510
$ws.Range("K92").Value = 7530
$ws.Range("M92").Value = -6282
$ws.Range("H122").Value = 723
$ws.Range("J122").Value = 800.8
$ws.Range("L122").Value = 7207.2
$ws.Range("N122").Value = -12107.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1649.4286
$ws.Range("I102").Value = 585.93335
$ws.Range("K102").Value = 585.93335
$ws.Range("M102").Value = 1036.06665
$ws.Range("H113").Value = 3044.125
$ws.Range("I113").Value = 1088.5
$ws.Range("K113").Value = 1088.5
$ws.Range("M113").Value = 1081.5
$ws.Range("H126").Value = 4785.4287
$ws.Range("I126").Value = 4499.3335
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 13498.0005
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -11028.0005
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2710.0625
$ws.Range("I132").Value = 2305.577
$ws.Range("K132").Value = 6916.731000000001
$ws.Range("M132").Value = -4386.731000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H55").Value = 499.70834
$ws.Range("J55").Value = 897.4286
$ws.Range("L55").Value = 897.4286
$ws.Range("N55").Value = -1243.4286
$ws.Range("H64").Value = 34944
$ws.Range("J64").Value = 34944
$ws.Range("L64").Value = 34944
$ws.Range("N64").Value = -35394
$ws.Range("H67").Value = 34944
$ws.Range("J67").Value = 34944
$ws.Range("L67").Value = 34944
$ws.Range("N67").Value = -36504
$ws.Range("H68").Value = 2631.125
$ws.Range("I68").Value = 2608.1667
$ws.Range("K68").Value = 2608.1667
$ws.Range("M68").Value = -1859.1667
$ws.Range("H71").Value = 2631.125
$ws.Range("I71").Value = 2608.1667
$ws.Range("K71").Value = 13040.8335
$ws.Range("M71").Value = -9296.833500000001
$ws.Range("H82").Value = 1486
$ws.Range("I82").Value = 1796.9231
$ws.Range("K82").Value = 1796.9231
$ws.Range("M82").Value = -1435.9231
$ws.Range("H85").Value = 1486
$ws.Range("I85").Value = 1796.9231
$ws.Range("K85").Value = 1796.9231
$ws.Range("M85").Value = -548.9231
$ws.Range("H132").Value = 5179.515
$ws.Range("I132").Value = 4754.7
$ws.Range("K132").Value = 14264.1
$ws.Range("M132").Value = -11734.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 890.125
$ws.Range("I107").Value = 469.1
$ws.Range("J107").Value = 1591.8334
$ws.Range("K107").Value = 1407.3
$ws.Range("L107").Value = 4775.5002
$ws.Range("M107").Value = 512.6999999999998
$ws.Range("N107").Value = -8615.5002
$ws.Range("H113").Value = 1444.2106
$ws.Range("I113").Value = 1139.5
$ws.Range("K113").Value = 3418.5
$ws.Range("M113").Value = -1248.5
$ws.Range("H122").Value = 1619.5
$ws.Range("I122").Value = 928
$ws.Range("J122").Value = 3002.5
$ws.Range("K122").Value = 2784
$ws.Range("L122").Value = 9007.5
$ws.Range("M122").Value = -334
$ws.Range("N122").Value = -13907.5
$ws.Range("H123").Value = 22663.334
$ws.Range("J123").Value = 22663.334
$ws.Range("L123").Value = 22663.334
$ws.Range("N123").Value = -32463.334
$ws.Range("H126").Value = 338999.66
$ws.Range("I126").Value = 1000000
$ws.Range("K126").Value = 3000000
$ws.Range("M126").Value = -2997530
$ws.Range("H132").Value = 1295.3125
$ws.Range("I132").Value = 1289.2727
$ws.Range("J132").Value = 1308.6
$ws.Range("K132").Value = 3867.8181
$ws.Range("L132").Value = 3925.8
$ws.Range("M132").Value = -1337.8181
$ws.Range("N132").Value = -8985.799999999999
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

